# "adding averages and more checks"
#
# Two kinds of change:
#  1) The dashboard's title/header styling is normalized onto a single
#     bold, white font (the old "bold, 14pt, default color" title-only
#     font is retired in favor of reusing the header's bold font, now
#     recolored white so it reads on the dark-blue header fill).
#  2) The training-expiry check was re-run on a later date, so the
#     "PERIOD TO EXPIRE" and "LAST UPDATE" columns for the first two
#     rows get refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# --- Title (A1): keep bold, drop the 14pt override, make it white ---
$titleFont = $ws.Range("A1").Font
$titleFont.Size = 11
$titleFont.Bold = $true
$titleFont.Color = 16777215

# --- Header row (A2:K2): bold white text on the existing dark fill ---
$headerFont = $ws.Range("A2:K2").Font
$headerFont.Bold = $true
$headerFont.Color = 16777215

# --- Refreshed data for rows 3 & 4 ---
# Leading "'" keeps these as literal text (matches the existing
# inline-string date cells) instead of Excel coercing them to date serials.
$ws.Range("H3").Value = -51
$ws.Range("I3").Value = "'16-Sep-2025"

$ws.Range("H4").Value = 300
$ws.Range("I4").Value = "'16-Sep-2025"
